$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2 through 442 all hold the same date
# serial value (45188 -> 2023-09-19). Update them all to 45189 (2023-09-20).
$ws.Range("C2:C442").Value = 45189
